$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "Anselmo-Gest. Int."
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("B4").Value = "Anselmo-Gest. Int."
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "Cleidson-Circuitos elétricos"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "[-, Valmir-Metrologia-1A, Valmir-Metrologia-1A, Valmir-Metrologia-1A]"
$ws.Range("D6").Value = "Cleidson-Circuitos elétricos"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "Cleidson-Circuitos elétricos"

# Row 7
$ws.Range("C7").Value = "[André Guimarães-Desenho técnico-1A, Joel L.-Tecnologia dos materiais-1A, Joel L.-Tecnologia dos materiais-1A, Valmir-Metrologia-1A]"
$ws.Range("F7").Value = "Cleidson-Circuitos elétricos"
